$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    93.50307287093942,
    93.634767339771727,
    93.678665496049163,
    93.371378402107112,
    93.064091308165047,
    93.371378402107112,
    93.678665496049163,
    93.854258121158907,
    93.678665496049163,
    93.766461808604035,
    93.854258121158907,
    93.590869183494291,
    93.415276558384548,
    93.415276558384548,
    93.459174714661984,
    94.600526777875331,
    94.600526777875331,
    94.644424934152767,
    93.634767339771727,
    93.590869183494291,
    93.634767339771727,
    94.117647058823522,
    93.810359964881471,
    93.590869183494291,
    93.810359964881471,
    93.50307287093942,
    93.634767339771727,
    93.942054433713778,
    94.424934152765587,
    94.249341527655844,
    93.50307287093942,
    93.371378402107112,
    93.590869183494291,
    93.678665496049163,
    93.634767339771727,
    84.899034240561903,
    93.590869183494291,
    93.722563652326613,
    93.898156277436357,
    93.634767339771727,
    93.546971027216856,
    93.546971027216856,
    93.50307287093942,
    93.546971027216856,
    93.371378402107112,
    93.766461808604035,
    84.240561896400351,
    93.810359964881471
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
